# Update the LR-pairs sheet with newly recomputed TPM-based expression values.
# The ligand (G/H) values are keyed by "Sending cluster", the receptor (M/N)
# values are keyed by "Target cluster". Every other numeric column in the
# sheet is derived from those two raw numbers, so we recompute the whole
# cascade for each data row:
#   H = G * (Ligand-expressing cells)
#   I = G / (sum of G over all sending clusters for this pair)
#   J = H / (sum of H over all sending clusters for this pair)
#   N = M * (Receptor-expressing cells)
#   O = M / (sum of M over all target clusters for this pair)
#   P = N / (sum of N over all target clusters for this pair)
#   Q = G * M   (edge average expression weight)
#   R = H * N   (edge total expression weight)
#   S = I * O   (edge average expression derived specificity)
#   T = J * P   (edge total expression derived specificity)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand average expression value (column G), per Sending cluster.
$ligandAvg = @{
    "ECs"           = 114.056483
    "FAPs"          = 251.7279513333333
    "MuSCs"         = 70.67310566666667
    "Resolving-Mac" = 1597.720744
}

# New receptor average expression value (column M), per Target cluster.
$receptorAvg = @{
    "ECs"           = 3.795192333333334
    "FAPs"          = 243.3763986666667
    "MuSCs"         = 29.801371
    "Resolving-Mac" = 86.47679266666667
}

$firstRow = 2
$lastRow = 17

# Sum of ligand/receptor averages across all clusters (shared denominator
# used by the specificity columns), since every row references the same
# four clusters for this single ligand-receptor pair.
$ligandTotal = 0
foreach ($key in $ligandAvg.Keys) { $ligandTotal += $ligandAvg[$key] }
$receptorTotal = 0
foreach ($key in $receptorAvg.Keys) { $receptorTotal += $receptorAvg[$key] }

$ligandTotalExpr = 0
foreach ($key in $ligandAvg.Keys) { $ligandTotalExpr += ($ligandAvg[$key] * 3) }
$receptorTotalExpr = 0
foreach ($key in $receptorAvg.Keys) { $receptorTotalExpr += ($receptorAvg[$key] * 3) }

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $sending = $ws.Cells.Item($row, 1).Value2
    $target = $ws.Cells.Item($row, 4).Value2

    $expressingCells = $ws.Cells.Item($row, 5).Value2  # Ligand-expressing cells (E)

    $G = $ligandAvg[$sending]
    $H = $G * $expressingCells
    $I = $G / $ligandTotal
    $J = $H / $ligandTotalExpr

    $M = $receptorAvg[$target]
    $N = $M * $expressingCells
    $O = $M / $receptorTotal
    $P = $N / $receptorTotalExpr

    $Q = $G * $M
    $R = $H * $N
    $S = $I * $O
    $T = $J * $P

    $ws.Cells.Item($row, 7).Value = $G    # G: Ligand average expression value
    $ws.Cells.Item($row, 8).Value = $H    # H: Ligand total expression value
    $ws.Cells.Item($row, 9).Value = $I    # I: Ligand derived specificity (avg)
    $ws.Cells.Item($row, 10).Value = $J   # J: Ligand derived specificity (total)

    $ws.Cells.Item($row, 13).Value = $M   # M: Receptor average expression value
    $ws.Cells.Item($row, 14).Value = $N   # N: Receptor total expression value
    $ws.Cells.Item($row, 15).Value = $O   # O: Receptor derived specificity (avg)
    $ws.Cells.Item($row, 16).Value = $P   # P: Receptor derived specificity (total)

    $ws.Cells.Item($row, 17).Value = $Q   # Q: Edge average expression weight
    $ws.Cells.Item($row, 18).Value = $R   # R: Edge total expression weight
    $ws.Cells.Item($row, 19).Value = $S   # S: Edge average expression derived specificity
    $ws.Cells.Item($row, 20).Value = $T   # T: Edge total expression derived specificity
}
